# Atualização dos dados: 26.12.2025 10:6
# Preenche a linha 13 (id 12) da planilha "quantidade_pontos" com os novos dados.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Valores simples (herdam o estilo padrao da coluna) ---
$ws.Range("B13").Value = 23
$ws.Range("C13").Value = 104495
$ws.Range("E13").Value = 96
$ws.Range("F13").Value = 1.5
$ws.Range("G13").Value = 1750
$ws.Range("I13").Value = 6980

# --- Formula (Valores multiplicados = Valores Totais * Multiplicador) ---
$ws.Range("D13").Formula = "=C13*F13"

# --- Tempo: copia a formatacao de hora da linha anterior e define o valor ---
$ws.Range("H12").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H13").Value = 0.056388888888888891
$excel.CutCopyMode = $false

# --- Data: copia a formatacao de data da linha anterior e define o valor ---
$ws.Range("L12").Copy()
$ws.Range("L13").PasteSpecial(-4122)
$ws.Range("L13").Value = 46017
$excel.CutCopyMode = $false

# --- Avatar / Tipo (texto) ---
$ws.Range("J13").Value = "Vampiro"
$ws.Range("K13").Value = "Desafio"

# Recalcula a planilha para atualizar a formula inserida
$excel.Calculate()

# Atualiza a celula selecionada, conforme o estado salvo do arquivo
$ws.Range("H25").Select()
